# Insert a new data row at row 8, pushing the existing rows 8..83 down to 9..84
# (old row 83 ends up as the new row 84). Then populate the new row 8 with its
# own values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(8).Insert()

$ws.Range('A8').Value = 5
$ws.Range('B8').Value = 'Macroferia Regional de Talca'
$ws.Range('C8').Value = 'Maule'
$ws.Range('D8').Value = 44882
$ws.Range('E8').Value = 7
$ws.Range('F8').Value = 100112040
$ws.Range('G8').Value = 'Cilantro'
$ws.Range('H8').Value = 'Sin especificar'
$ws.Range('I8').Value = 'Primera'
$ws.Range('J8').Value = 150
$ws.Range('K8').Value = 7000
$ws.Range('L8').Value = 7000
$ws.Range('M8').Value = 7000
$ws.Range('N8').Value = '$/caja 36 atados'
$ws.Range('O8').Value = 'Región del Maule'
$ws.Range('P8').Value = 194
$ws.Range('Q8').Value = 36
$ws.Range('R8').Value = 'Hortaliza'
